# Applies the "Holden scheme" update to UniformA-HW10.xlsx
#  - removes the now-unused duplicate columns X:AG (row 1 numbering / row 2 duplicate headers)
#  - reorders the [h,k,l] plane headers in row 2 (C2:W2)
#  - renames rows 16:19 from the HexGrid-* entries to the new Holden2.5/5/10/15 entries
#  - appends 4 new rows (20:23) carrying the HexGrid-* entries that got displaced

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the stray duplicate block in columns X:AG (rows 1-19) so the sheet
#    shrinks back down to a A1:W.. extent, matching the new dimension.
$ws.Range("X1:AG19").Delete()

# 2) Row 2 headers (C2:W2): the [h,k,l] plane columns are reshuffled; the
#    "Pairs" columns (N2:W2) keep their original order/values.
$row2Labels = @(
    "[4, 2, 0]", "[4, 0, 0]", "[2, 0, 0]", "[2, 2, 0]", "[3, 3, 3]",
    "[1, 1, 1]", "[2, 2, 2]", "[3, 3, 1]", "[3, 1, 1]", "[4, 2, 2]", "[5, 1, 1]",
    "1Pair-A", "1Pair-B", "2Pairs-A", "2Pairs-B",
    "3Pairs-A", "3Pairs-B", "3Pairs-C", "4Pairs", "5A4F", "MaxUnique"
)
for ($i = 0; $i -lt $row2Labels.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $row2Labels[$i]
}

# 3) Rows 16:19 used to be the HexGrid-* scans; they now become the new
#    Holden2.5 / Holden5 / Holden10 / Holden15 scans (same A/C:W data shape).
$holdenLabels = @("Holden2.5", "Holden5", "Holden10", "Holden15")
for ($i = 0; $i -lt $holdenLabels.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = $holdenLabels[$i]
}

# 4) Append 4 new rows (20:23) holding the displaced HexGrid-* scans. Clone
#    the formatting of the last existing data row (19) so the new A/B cells
#    pick up the same bold/border/centred style used throughout column A
#    and the row header column B.
$hexLabels = @(
    "HexGrid-90degTilt2.5degRes",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt10degRes",
    "HexGrid-90degTilt15degRes"
)
for ($i = 0; $i -lt $hexLabels.Length; $i++) {
    $r = 20 + $i
    $ws.Range("A19:W19").Copy()
    $ws.Range("A" + $r + ":W" + $r).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = 18 + $i
    $ws.Cells.Item($r, 2).Value = $hexLabels[$i]
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
